$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# Append the new log row (row 10) mirroring the existing column layout.
$ws.Range("A10").Value = "Geen onderwerp"
$ws.Range("B10").Value = "onbekend"
$ws.Range("D10").Value = "Onbekend"
$ws.Range("F10").Value = "2025-08-18 21:20:52"
$ws.Range("G10").Value = "Nee"
$ws.Range("H10").Value = "Ja"
$ws.Range("I10").Value = "Nee"
$ws.Range("J10").Value = "Nee"

# Extend the conditional-formatting ranges so row 10 is covered too.
$ws.Range("D2:D9").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D10"))
$ws.Range("G2:G9").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G10"))
$ws.Range("H2:H9").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H10"))
$ws.Range("I2:I9").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I10"))
$ws.Range("J2:J9").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J10"))

# Update the Dashboard summary count for the "Onbekend" category (4 -> 5).
$dash.Range("B2").Value = 5
